$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with two new values continuing the sequence
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy the bold/bordered header formatting from O1 onto the new P1:Q1 cells
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# For each data row (2..25) swap columns I/K and M/O, and add new columns P and Q
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q = 2
}
